$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 3: "Current " + "Job Openings:" -> "Current Job Openings:"
# Use an absolute Characters() sub-range (rather than Paragraphs(3) directly)
# so the merge happens even though the resulting text is unchanged.
$para3 = $tr.Paragraphs(3)
$para3Range = $tr.Characters($para3.Start, $para3.Length)
$para3Range.Text = "Current Job Openings:"

# Paragraph 6: "...WebKit" + " " + "Development, Quality and automation"
# Only merge the last two runs (" " and "Development, Quality and automation"),
# leaving "Computer Scientist - " and "WebKit" runs untouched.
$para6 = $tr.Paragraphs(6)
$mergeStart = $para6.Start + 21 + 6
$mergeLen = 1 + 35
$mergeRange = $tr.Characters($mergeStart, $mergeLen)
$mergeRange.Text = " Development, Quality and automation"
